{"js": "const replacements = [\n  [\"721\u00f76=120, 1\", \"553\u00f72=276, 1\"],\n  [\"545\u00f76=90, 5\", \"466\u00f73=155, 1\"],\n  [\"719\u00f79=79, 8\", \"682\u00f73=227, 1\"],\n  [\"911\u00f75=182, 1\", \"715\u00f74=178, 3\"],\n  [\"616\u00f75=123, 1\", \"834\u00f76=139, 0\"],\n  [\"951\u00f74=237, 3\", \"763\u00f72=381, 1\"],\n  [\"741\u00f77=105, 6\", \"577\u00f77=82, 3\"],\n  [\"896\u00f76=149, 2\", \"889\u00f74=222, 1\"],\n  [\"218\u00f73=72, 2\", \"774\u00f73=258, 0\"],\n  [\"875\u00f77=125, 0\", \"871\u00f76=145, 1\"],\n  [\"800\u00f73=266, 2\", \"764\u00f77=109, 1\"],\n  [\"937\u00f73=312, 1\", \"843\u00f75=168, 3\"],\n  [\"706\u00f74=176, 2\", \"309\u00f73=103, 0\"],\n  [\"917\u00f78=114, 5\", \"701\u00f79=77, 8\"],\n  [\"450\u00f74=112, 2\", \"122\u00f75=24, 2\"],\n  [\"868\u00f74=217, 0\", \"651\u00f75=130, 1\"],\n  [\"285\u00f76=47, 3\", \"779\u00f75=155, 4\"],\n  [\"697\u00f79=77, 4\", \"261\u00f76=43, 3\"],\n  [\"308\u00f75=61, 3\", \"553\u00f76=92, 1\"],\n  [\"253\u00f79=28, 1\", \"372\u00f79=41, 3\"],\n  [\"892\u00f76=148, 4\", \"265\u00f72=132, 1\"],\n  [\"943\u00f75=188, 3\", \"689\u00f78=86, 1\"],\n  [\"920\u00f75=184, 0\", \"158\u00f74=39, 2\"],\n  [\"513\u00f77=73, 2\", \"962\u00f73=320, 2\"],\n  [\"715\u00f73=238, 1\", \"303\u00f76=50, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"721\u00f76=120, 1\", \"553\u00f72=276, 1\"),\n  @(\"545\u00f76=90, 5\", \"466\u00f73=155, 1\"),\n  @(\"719\u00f79=79, 8\", \"682\u00f73=227, 1\"),\n  @(\"911\u00f75=182, 1\", \"715\u00f74=178, 3\"),\n  @(\"616\u00f75=123, 1\", \"834\u00f76=139, 0\"),\n  @(\"951\u00f74=237, 3\", \"763\u00f72=381, 1\"),\n  @(\"741\u00f77=105, 6\", \"577\u00f77=82, 3\"),\n  @(\"896\u00f76=149, 2\", \"889\u00f74=222, 1\"),\n  @(\"218\u00f73=72, 2\", \"774\u00f73=258, 0\"),\n  @(\"875\u00f77=125, 0\", \"871\u00f76=145, 1\"),\n  @(\"800\u00f73=266, 2\", \"764\u00f77=109, 1\"),\n  @(\"937\u00f73=312, 1\", \"843\u00f75=168, 3\"),\n  @(\"706\u00f74=176, 2\", \"309\u00f73=103, 0\"),\n  @(\"917\u00f78=114, 5\", \"701\u00f79=77, 8\"),\n  @(\"450\u00f74=112, 2\", \"122\u00f75=24, 2\"),\n  @(\"868\u00f74=217, 0\", \"651\u00f75=130, 1\"),\n  @(\"285\u00f76=47, 3\", \"779\u00f75=155, 4\"),\n  @(\"697\u00f79=77, 4\", \"261\u00f76=43, 3\"),\n  @(\"308\u00f75=61, 3\", \"553\u00f76=92, 1\"),\n  @(\"253\u00f79=28, 1\", \"372\u00f79=41, 3\"),\n  @(\"892\u00f76=148, 4\", \"265\u00f72=132, 1\"),\n  @(\"943\u00f75=188, 3\", \"689\u00f78=86, 1\"),\n  @(\"920\u00f75=184, 0\", \"158\u00f74=39, 2\"),\n  @(\"513\u00f77=73, 2\", \"962\u00f73=320, 2\"),\n  @(\"715\u00f73=238, 1\", \"303\u00f76=50, 3\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    Write-Output \"Not found: $oldText\"\n  }\n}"}
